$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it sat around the "Comparar
#    aproximaciones" image earlier in the document). We'll re-add it at its
#    new location later.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Count -gt 0) {
    $d.Bookmarks.Item("_GoBack").Delete()
} else {
    try {
        $d.Bookmarks.Item("_GoBack").Delete()
    } catch {
    }
}

# ---------------------------------------------------------------------------
# 2) Replace the single sentence
#      "Se cerrará la aplicación. Todos los datos utilizados se borrarán."
#    with the new, longer explanation made up of several runs (some
#    italicised) describing the "Nueva sesión" / "Finalizar" choice.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Se cerrará la aplicación. Todos los datos utilizados se borrarán.")

if ($found) {
    # First chunk re-uses the found range itself.
    $rng.Text = "Se "

    $rng.Collapse(0)
    $rng.InsertAfter("le preguntará si quiere iniciar una ")

    $rng.Collapse(0)
    $rng.InsertAfter("Nueva sesión")
    $rng.Italic = 1

    $rng.Collapse(0)
    $rng.InsertAfter(" ")
    $rng.Italic = 1

    $rng.Collapse(0)
    $rng.InsertAfter("(borrar todos los datos ingresados anteriormente e ingresar nuevamente para trabajar con nuevos datos) o ")
    $rng.Italic = 0

    $rng.Collapse(0)
    $rng.InsertAfter("Finalizar")
    $rng.Italic = 1

    $rng.Collapse(0)
    $rng.InsertAfter(" ")
    $rng.Italic = 0
    $rng.Font.Size = 28

    $rng.Collapse(0)
    $rng.InsertAfter("(")
    $rng.Font.Size = 24

    $rng.Collapse(0)
    $rng.InsertAfter("cerrar el programa).")

    # -----------------------------------------------------------------------
    # 3) Re-insert the "_GoBack" bookmark right after the new text.
    # -----------------------------------------------------------------------
    $rng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng)
}

# ---------------------------------------------------------------------------
# 4) The footer's cached PAGE field result changes from 8 to 7.
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$pageFound = $ftr.Range.Find.Execute("8")
if ($pageFound) {
    $ftr.Range.Characters.Item(1).Text = "7"
}
